$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[[0.44889334]`n [0.441223  ]`n [0.45822978]]"
$ws.Range("C2").Value = "[[0.20702952]`n [0.17418403]`n [0.1492001 ]]"

# Setting a value containing embedded line breaks triggers Excel's automatic
# row-height expansion (even without WrapText). Re-running AutoFit restores
# the row to its natural/default height so no spurious ht/customHeight
# attributes are introduced, keeping the change scoped to the two cell values.
$ws.Rows.Item(2).AutoFit()
